$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Data row (admin already exists as shared string index 0; manager is new)
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "FailMsg"

$ws.Range("C2").Value = "Home Page is not displayed"

# Auto-fit the columns that hold the longer text (matches the authored column widths)
$ws.Columns("A").AutoFit()
$ws.Columns("C").AutoFit()

# Match the resulting selection from the diff
$ws.Range("A3").Select()
